$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsInfo = $wb.Worksheets.Item("ThongTinBang")

# ---------------------------------------------------------------------------
# Sheet "Data": update row 2 values, then drop the old sample rows 3-6.
# ---------------------------------------------------------------------------
$wsData.Range("A2").Value = "danh muc goc"
$wsData.Range("B2").Value = "dmgoc"
$wsData.Range("D2").Value = 0
$wsData.Range("E2").Value = 123

$wsData.Rows("3:6").Delete() | Out-Null

# Shrink Table1 down to the single remaining data row.
$tbl1 = $wsData.ListObjects.Item("Table1")
$tbl1.Resize($wsData.Range("A1:E2")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "ThongTinBang": fix the accented typo and append the new categories.
# ---------------------------------------------------------------------------
$wsInfo.Range("A3").Value = "danh muc 1"

$wsInfo.Range("A4").Value = "danh muc 2"
$wsInfo.Range("B4").Value = "dm2"
$wsInfo.Range("A5").Value = "danh muc 3"
$wsInfo.Range("B5").Value = "dm3"
$wsInfo.Range("A6").Value = "danh muc 4"
$wsInfo.Range("B6").Value = "dm4"

# Grow Table2 (TenDanhMucCha / MaDanhMuc) to cover the new rows.
$tbl2 = $wsInfo.ListObjects.Item("Table2")
$tbl2.Resize($wsInfo.Range("A1:B6")) | Out-Null

# ---------------------------------------------------------------------------
# Selections: ThongTinBang!B4 selected (inactive tab), Data!F4 selected and
# the active tab - matches the final saved view state.
# ---------------------------------------------------------------------------
$wsInfo.Activate() | Out-Null
$wsInfo.Range("B4").Select() | Out-Null

$wsData.Activate() | Out-Null
$wsData.Range("F4").Select() | Out-Null
